# edit.ps1 - apply the "Add protype of project-presentation" commit to idea.docx
#
# The target diff does two independent things:
#   1. Tags the title/subtitle runs on the cover page with an explicit
#      en-GB language and splices a new "STORY BASED " run into the
#      subtitle ("A ROUGELIKE BULLETHELL RPG" -> "A STORY BASED
#      ROUGELIKE BULLETHELL RPG").
#   2. Refreshes the table of contents: new TOC-bookmark ids, a new
#      "Stile" TOC entry, "Handlung" -> "Handlung (Exzerpt)" in the TOC
#      text (matching the heading, whose bookmark now spans the whole
#      "Handlung (Exzerpt)" run sequence), and the _GoBack bookmark
#      hopping from the end of the document to just after the TOC
#      (where the cursor was when the file was last saved).
#
# Each touched paragraph is replaced wholesale with Range.InsertXML,
# using the paragraph's own original OOXML (captured verbatim) with only
# the minimal diff-described edits applied, so everything this host
# doesn't need to change (rsids, run splits, formatting, field codes...)
# survives byte-for-byte.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaXml($index, $xml) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) Title paragraph "CODENAME: HERO" -> add <w:lang w:val="en-GB"/>
# ---------------------------------------------------------------------
$xml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Titel"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="DejaVu Serif" w:hAnsi="DejaVu Serif"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="DejaVu Serif" w:hAnsi="DejaVu Serif"/><w:lang w:val="en-GB"/></w:rPr><w:t>CODENAME: HERO</w:t></w:r></w:p>'
Set-ParaXml 4 $xml

# ---------------------------------------------------------------------
# 2) Subtitle paragraph "A ROUGELIKE BULLETHELL RPG" -> tag runs with
#    en-GB and splice in a new "STORY BASED " run.
# ---------------------------------------------------------------------
$xml = '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="Untertitel"/><w:jc w:val="center"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">STORY BASED </w:t></w:r><w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>ROUGELIKE BULLETHELL RPG</w:t></w:r></w:p>'
Set-ParaXml 5 $xml

# ---------------------------------------------------------------------
# 3) The blank paragraph right after the subtitle (spacing after=0,
#    centered, DejaVu Serif) -> add <w:lang w:val="en-GB"/> too.
# ---------------------------------------------------------------------
$xml = '<w:p ' + $wNs + '><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="DejaVu Serif" w:hAnsi="DejaVu Serif"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p>'
Set-ParaXml 6 $xml

# ---------------------------------------------------------------------
# TOC / heading-bookmark refresh. Walk bottom-up so paragraph indices
# for not-yet-touched paragraphs never shift under us.
# ---------------------------------------------------------------------

# 51) "Titan Souls" list item -> drop the stale trailing _GoBack bookmark
$xml = '<w:p w:rsidR="0011149D" w:rsidRPr="00E13968" w:rsidRDefault="00871A93" w:rsidP="0011149D"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Titan Souls</w:t></w:r></w:p>'
Set-ParaXml 51 $xml

# 41) "Stile" heading -> wrap it in the new _Toc462568097 bookmark
$xml = '<w:p w:rsidR="00E13968" w:rsidRDefault="00E13968" w:rsidP="00871A93"><w:pPr><w:pStyle w:val="berschrift1"/><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr></w:pPr><w:bookmarkStart w:id="3" w:name="_Toc462568097"/><w:r w:rsidRPr="00E13968"><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr><w:t>Stile</w:t></w:r><w:bookmarkEnd w:id="3"/></w:p>'
Set-ParaXml 41 $xml

# 38) "Handlung (Exzerpt)" heading -> renumber bookmark id/name and move
#     the bookmarkEnd out to behind "(Exzerpt)" (it used to close right
#     after "Handlung").
$xml = '<w:p w:rsidR="009A7ADB" w:rsidRDefault="00B665D8" w:rsidP="00B665D8"><w:pPr><w:pStyle w:val="berschrift1"/><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr></w:pPr><w:bookmarkStart w:id="2" w:name="_Toc462568096"/><w:r w:rsidRPr="00B665D8"><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr><w:t>Handlung</w:t></w:r><w:r w:rsidR="008B1CDE"><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="008B1CDE" w:rsidRPr="00E52621"><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/><w:sz w:val="28"/></w:rPr><w:t>(Exzerpt)</w:t></w:r><w:bookmarkEnd w:id="2"/></w:p>'
Set-ParaXml 38 $xml

# 36) "Beschreibung" heading -> renumber bookmark id/name
$xml = '<w:p w:rsidR="00BE0B1B" w:rsidRPr="0012744A" w:rsidRDefault="00181985" w:rsidP="00181985"><w:pPr><w:pStyle w:val="berschrift1"/><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr></w:pPr><w:bookmarkStart w:id="1" w:name="_Toc462568095"/><w:r w:rsidRPr="0012744A"><w:rPr><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:color w:val="auto"/></w:rPr><w:lastRenderedPageBreak/><w:t>Beschreibung</w:t></w:r><w:bookmarkEnd w:id="1"/></w:p>'
Set-ParaXml 36 $xml

# 35) page-break paragraph right before the "Beschreibung" heading ->
#     this is where the cursor sat when the file was saved, so _GoBack
#     moves here from the end of the document.
$xml = '<w:p w:rsidR="00181985" w:rsidRDefault="00181985" w:rsidP="00181985"><w:r><w:br w:type="page"/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-ParaXml 35 $xml

# 34) the field-end paragraph right after the "Handlung" TOC entry ->
#     split into the new "Stile" TOC entry + the (moved) field-end
#     paragraph.
$xml = '<w:p><w:pPr><w:pStyle w:val="Verzeichnis1"/><w:tabs><w:tab w:val="right" w:pos="9062"/></w:tabs><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="de-AT"/></w:rPr></w:pPr><w:hyperlink w:anchor="_Toc462568097" w:history="1"><w:r><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:noProof/></w:rPr><w:t>Stile</w:t></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:instrText xml:space="preserve"> PAGEREF _Toc462568097 \h </w:instrText></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:hyperlink></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p>'
Set-ParaXml 34 $xml

# 33) "Handlung" TOC entry -> renumber PAGEREF/anchor, text gains
#     " (Exzerpt)" to mirror the heading it now points at.
$xml = '<w:p w:rsidR="00C140C4" w:rsidRDefault="00EA7C45"><w:pPr><w:pStyle w:val="Verzeichnis1"/><w:tabs><w:tab w:val="right" w:pos="9062"/></w:tabs><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="de-AT"/></w:rPr></w:pPr><w:hyperlink w:anchor="_Toc462568096" w:history="1"><w:r w:rsidR="00C140C4" w:rsidRPr="00507F6C"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:noProof/></w:rPr><w:t>Handlung (Exzerpt)</w:t></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:tab/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:instrText xml:space="preserve"> PAGEREF _Toc462568096 \h </w:instrText></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:hyperlink></w:p>'
Set-ParaXml 33 $xml

# 32) "Beschreibung" TOC entry (also carries the TOC field begin /
#     instrText) -> renumber PAGEREF/anchor.
$xml = '<w:p w:rsidR="00C140C4" w:rsidRDefault="00181985"><w:pPr><w:pStyle w:val="Verzeichnis1"/><w:tabs><w:tab w:val="right" w:pos="9062"/></w:tabs><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorEastAsia" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:b w:val="0"/><w:bCs w:val="0"/><w:caps w:val="0"/><w:noProof/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:eastAsia="de-AT"/></w:rPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> TOC \o "1-3" \h \z </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:hyperlink w:anchor="_Toc462568095" w:history="1"><w:r w:rsidR="00C140C4" w:rsidRPr="00507F6C"><w:rPr><w:rStyle w:val="Hyperlink"/><w:rFonts w:ascii="DejaVu Sans Mono" w:hAnsi="DejaVu Sans Mono" w:cs="DejaVu Sans Mono"/><w:noProof/></w:rPr><w:t>Beschreibung</w:t></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:tab/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:instrText xml:space="preserve"> PAGEREF _Toc462568095 \h </w:instrText></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:t>2</w:t></w:r><w:r w:rsidR="00C140C4"><w:rPr><w:noProof/><w:webHidden/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:hyperlink></w:p>'
Set-ParaXml 32 $xml

Write-Output "all steps done"
